# Second commit: set A1 text, change sheet font to Arial, resize column A,
# and leave the selection on C6 (matching the saved sheetView selection).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A1 gets the "name : marwa " text (stored as a shared string)
$ws.Range("A1").Value = "name : marwa "

# Default font changes from Calibri to Arial
$ws.Cells.Font.Name = "Arial"

# Column A is widened to ~16.6 characters
$ws.Columns.Item(1).ColumnWidth = 15.83

# Final selection/active cell is C6
[void]$ws.Range("C6").Select()
